{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// First paragraph: \"This is a Microsoft word document.\"\nconst firstPara = paragraphs.items[0];\n\n// Append two trailing spaces to the existing text (stays in the same run).\nconst tail = firstPara.getRange(\"End\");\ntail.insertText(\"  \", \"End\");\nawait context.sync();\n\n// Insert the new, differently-colored annotation as its own run right after.\nconst endOfPara = firstPara.getRange(\"End\");\nconst newRun = endOfPara.insertText(\n  \"(This is a change \\u2013 Version for branch alternate)\",\n  \"End\"\n);\nnewRun.font.color = \"#C00000\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# First paragraph: \"This is a Microsoft word document.\"\n$firstPara = $d.Paragraphs.Item(1).Range\n\n# Append two trailing spaces; InsertAfter lands the text before the\n# paragraph mark, so this stays inside the same (first) run.\n$firstPara.InsertAfter(\"  \")\n\n# Append the new annotation text right after the padding, still inside\n# paragraph 1.\n$newText = \"(This is a change \u2013 Version for branch alternate)\"\n$firstPara.InsertAfter($newText)\n\n# Locate the just-inserted annotation and recolor just that span, which\n# splits it into its own run distinct from the original sentence.\n$annotationRange = $d.Content\n$annotationRange.Find.ClearFormatting()\n$annotationRange.Find.Text = $newText\n$annotationRange.Find.Execute() | Out-Null\n$annotationRange.Font.Color = 192\n"}
